$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split each `"Fraction e(8, 2); //"` run into three runs so the
#    literal "2" becomes its own run, then retype it as "4":
#       "Fraction e(8,   +   4   +   ); //"
#    (matches the two code-listing screenshots in the doc; sizes 19/20)
# ---------------------------------------------------------------------
$needle = '"Fraction e(8, 2); //"'
$prefixLen = '"Fraction e(8, '.Length

$search = $d.Content
$found = $search.Find.Execute($needle)
while ($found) {
    $matchStart = $search.Start
    $digitStart = $matchStart + $prefixLen
    $digitEnd = $digitStart + 1

    # Replace the "2" character in place (keeps its existing formatting).
    $digitRange = $d.Range($digitStart, $digitEnd)
    $digitRange.Text = "4"

    # Force Word to materialize a dedicated run for the new character
    # (toggling a character property off again splits the run without
    # leaving any stray formatting behind).
    $digitRange2 = $d.Range($digitStart, $digitStart + 1)
    $digitRange2.Font.Bold = $true
    $digitRange2.Font.Bold = $false

    $search = $d.Range($digitStart + 1, $d.Content.End)
    $found = $search.Find.Execute($needle)
}

# ---------------------------------------------------------------------
# 2) Plain-text occurrences: `Fraction e(8, 2); // e = 2`
#    -> `Fraction e(8, 4); // e = 2` (no run split here).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Fraction e(8, 2); // e = 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fraction e(8, 4); // e = 2", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Move the `_GoBack` bookmark from the top of the document (right
#    before "Part C - Encapsulation") down to just after the "00" run
#    inside "...proflastname/submit 2005_w5_home..." (collapsed,
#    zero-length bookmark). Adding a bookmark with the same name moves
#    it, removing the old one automatically.
# ---------------------------------------------------------------------
$anchorNeedle = "proflastname/submit 200_w5_home"
$anchor = $d.Content
$anchor.Find.Execute($anchorNeedle) | Out-Null
$bookmarkPos = $anchor.Start + "proflastname/submit 200".Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
